# Apply updated cryptocurrency market data to the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin name / link swap for rows 23-24 (Uniswap <-> Avalanche reordered upstream)
$ws.Range('B23').Value = 'Avalanche'
$ws.Range('C23').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('B24').Value = 'Uniswap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'

# Updated Price (column D) and Volume(1h) (column E) values
# Each entry: cell reference, new value, whether the text looks like a plain
# number (requires forcing Text format so Excel does not convert it)
$updates = @(
    @('D2', '34.507.58', $false),
    @('E2', '  +13.74%  ', $false),
    @('D3', '1.830.18', $false),
    @('E3', '  +9.09%  ', $false),
    @('D4', '0.997', $true),
    @('E4', '  -0.11%  ', $false),
    @('D5', '235.30', $true),
    @('E5', '  +7.21%  ', $false),
    @('E6', '  +5.74%  ', $false),
    @('D7', '0.999', $true),
    @('E7', '  +0.13%  ', $false),
    @('D8', '32.10', $true),
    @('E8', '  +7.60%  ', $false),
    @('D9', '46.32', $true),
    @('E9', '  +5.32%  ', $false),
    @('D10', '0.287', $true),
    @('E10', '  +8.57%  ', $false),
    @('D11', '0.0688', $true),
    @('E11', '  +11.17%  ', $false),
    @('D12', '0.0930', $true),
    @('E12', '  +2.73%  ', $false),
    @('D13', '2.077.46', $false),
    @('E13', '  +8.34%  ', $false),
    @('D14', '1.820.29', $false),
    @('E14', '  +8.58%  ', $false),
    @('D15', '0.650', $true),
    @('E15', '  +5.09%  ', $false),
    @('D16', '34.517.76', $false),
    @('E16', '  +13.79%  ', $false),
    @('D17', '10.36', $true),
    @('E17', '  -3.49%  ', $false),
    @('D18', '4.39', $true),
    @('E18', '  +9.45%  ', $false),
    @('D19', '71.67', $true),
    @('E19', '  +8.90%  ', $false),
    @('D20', '266.13', $true),
    @('E20', '  +7.07%  ', $false),
    @('D21', '0.0₃0764', $false),
    @('E21', '  +6.21%  ', $false),
    @('D22', '0.995', $true),
    @('E22', '  -0.29%  ', $false),
    @('D23', '10.65', $true),
    @('E23', '  +5.26%  ', $false),
    @('D24', '4.46', $true),
    @('E24', '  +3.37%  ', $false),
    @('D25', '2.18', $true),
    @('E25', '  -1.76%  ', $false),
    @('D26', '162.78', $true),
    @('E26', '  +2.35%  ', $false),
    @('D27', '17.16', $true),
    @('E27', '  +8.02%  ', $false),
    @('D28', '0.118', $true),
    @('E28', '  +5.80%  ', $false),
    @('D29', '7.21', $true),
    @('E29', '  +6.69%  ', $false),
    @('E30', '  -0.43%  ', $false),
    @('D31', '3.89', $true),
    @('E31', '  +11.55%  ', $false),
    @('D32', '0.0521', $true),
    @('E32', '  +4.18%  ', $false),
    @('D33', '1.22', $true),
    @('E33', '  +7.13%  ', $false),
    @('E34', '  +9.40%  ', $false),
    @('D35', '1.615.27', $false),
    @('E35', '  +8.85%  ', $false),
    @('D36', '1.88', $true),
    @('E36', '  +8.50%  ', $false),
    @('D37', '89.88', $true),
    @('E37', '  +12.84%  ', $false),
    @('E38', '  +3.41%  ', $false),
    @('D39', '0.635', $true),
    @('E39', '  +7.67%  ', $false),
    @('D40', '0.0190', $true),
    @('E40', '  +5.75%  ', $false),
    @('D41', '2.86', $true),
    @('E41', '  +5.55%  ', $false),
    @('D42', '0.933', $true),
    @('E42', '  +8.62%  ', $false),
    @('E43', '  +2.56%  ', $false),
    @('D44', '2.18', $true),
    @('E44', '  +7.59%  ', $false),
    @('D45', '0.0521', $true),
    @('E45', '  +3.08%  ', $false),
    @('E46', '  +3.64%  ', $false),
    @('D47', '1.967.96', $false),
    @('E47', '  +8.39%  ', $false),
    @('D48', '54.63', $true),
    @('E48', '  +3.31%  ', $false),
    @('D49', '5.83', $true),
    @('E49', '  +7.02%  ', $false),
    @('D50', '1.00', $true),
    @('E50', '  +0.13%  ', $false),
    @('D51', '11.65', $true),
    @('E51', '  +24.66%  ', $false)
)

foreach ($u in $updates) {
    $cellRef = $u[0]
    $newValue = $u[1]
    $looksNumeric = $u[2]
    $rng = $ws.Range($cellRef)
    if ($looksNumeric) {
        # Force text storage so values like '1.00' or '4.46' are not
        # reinterpreted as numbers, then restore default styling so no
        # residual cell-format attribute is left behind.
        $rng.NumberFormat = '@'
        $rng.Value = $newValue
        $rng.Style = 'Normal'
    } else {
        $rng.Value = $newValue
    }
}
